$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3").Value = "new branch line added here."
$ws.Range("E3").Select()
